$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("2023-02-08", $true, $false, $false, $false, $false, $true, 1, $false, "2023-02-14", 2)
